# TC31_Canine_Filter_Breed-IrishSettr.xlsx
# Replace the old per-tab "StatQuery" Cypher query (column C, rows 2-4) with
# the new Programs/Studies/Cases/Samples/Case Files/Study Files query, fix
# the sheet zoom back to 100%, and move the active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Irish Setter']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Zoom back to 100% (was saved at 70%) and move the selection to B4.
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("B4").Select()
